$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as text so values like
# "216.95" or "0.520" are not auto-converted to numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.878.06"
$ws.Range("E2").Value = "  -0.13%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.650.50"
$ws.Range("E3").Value = "  +1.47%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "216.95"
$ws.Range("E5").Value = "  +1.22%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.520"
$ws.Range("E6").Value = "  -0.32%  "

# Row 7 - USDC
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - Solana
$ws.Range("D8").Value = "29.03"
$ws.Range("E8").Value = "  -2.28%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.262"
$ws.Range("E9").Value = "  +1.34%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.0610"
$ws.Range("E10").Value = "  -0.17%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0899"
$ws.Range("E11").Value = "  -1.72%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.875.95"
$ws.Range("E12").Value = "  +0.92%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.636.06"
$ws.Range("E13").Value = "  +0.55%  "

# Row 14 - now Chainlink (was Polygon)
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "10.02"
$ws.Range("E14").Value = "  +12.41%  "

# Row 15 - now Polygon (was Chainlink)
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.596"
$ws.Range("E15").Value = "  +4.06%  "

# Row 16 - Polkadot
$ws.Range("D16").Value = "3.92"
$ws.Range("E16").Value = "  +0.48%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "29.934.13"
$ws.Range("E17").Value = "  -0.10%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "64.49"
$ws.Range("E18").Value = "  -0.40%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "239.39"
$ws.Range("E19").Value = "  -1.96%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0707"
$ws.Range("E20").Value = "  +0.07%  "

# Row 21 - Dai
$ws.Range("D21").Value = "0.997"
$ws.Range("E21").Value = "  -0.06%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "9.96"
$ws.Range("E22").Value = "  +3.23%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "4.15"
$ws.Range("E23").Value = "  +0.17%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "2.18"
$ws.Range("E24").Value = "  +2.36%  "

# Row 25 - Monero
$ws.Range("D25").Value = "157.73"
$ws.Range("E25").Value = "  +0.11%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "15.70"
$ws.Range("E26").Value = "  +0.22%  "

# Row 27 - Stellar (D unchanged)
$ws.Range("E27").Value = "  -0.80%  "

# Row 28 - Cosmos
$ws.Range("D28").Value = "6.71"
$ws.Range("E28").Value = "  +1.48%  "

# Row 29 - BinanceUSD
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  -0.16%  "

# Row 30 - Hedera (D unchanged)
$ws.Range("E30").Value = "  +1.41%  "

# Row 31 - PancakeSwap (D unchanged)
$ws.Range("E31").Value = "  -0.70%  "

# Row 32 - Filecoin (D unchanged)
$ws.Range("E32").Value = "  +1.22%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "3.21"
$ws.Range("E33").Value = "  -0.90%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.428.54"
$ws.Range("E34").Value = "  +0.46%  "

# Row 35 - LidoDAOToken
$ws.Range("D35").Value = "1.69"
$ws.Range("E35").Value = "  +3.36%  "

# Row 36 - TrustWalletToken (D unchanged)
$ws.Range("E36").Value = "  -0.97%  "

# Row 37 - now MXToken (was VeChain)
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "2.70"
$ws.Range("E37").Value = "  -6.05%  "

# Row 38 - now VeChain (was MXToken)
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.0174"
$ws.Range("E38").Value = "  +1.93%  "

# Row 39 - HuobiToken
$ws.Range("D39").Value = "2.30"
$ws.Range("E39").Value = "  +0.15%  "

# Row 40 - ImmutableX
$ws.Range("D40").Value = "0.574"
$ws.Range("E40").Value = "  +2.43%  "

# Row 41 - Aave
$ws.Range("D41").Value = "77.39"
$ws.Range("E41").Value = "  +11.25%  "

# Row 42 - ARBITRUM
$ws.Range("D42").Value = "0.841"
$ws.Range("E42").Value = "  +0.86%  "

# Row 43 - Kaspa
$ws.Range("D43").Value = "0.0502"
$ws.Range("E43").Value = "  +0.05%  "

# Row 44 - RenderToken
$ws.Range("D44").Value = "1.95"
$ws.Range("E44").Value = "  -2.41%  "

# Row 45 - PaxDollar
$ws.Range("D45").Value = "0.996"
$ws.Range("E45").Value = "  -0.05%  "

# Row 46 - WEMIXToken (D unchanged)
$ws.Range("E46").Value = "  -2.42%  "

# Row 47 - BitcoinSV
$ws.Range("D47").Value = "50.53"
$ws.Range("E47").Value = "  -7.31%  "

# Row 48 - RocketPoolETH
$ws.Range("D48").Value = "1.784.44"
$ws.Range("E48").Value = "  +1.01%  "

# Row 49 - FraxShare
$ws.Range("D49").Value = "5.35"
$ws.Range("E49").Value = "  -0.97%  "

# Row 50 - Quant
$ws.Range("D50").Value = "93.97"
$ws.Range("E50").Value = "  +5.59%  "

# Row 51 - BabyDogeCoin
$ws.Range("D51").Value = "0.0₆0108"
$ws.Range("E51").Value = "  -0.79%  "

# Restore the default (no explicit number format) style so the written
# cells match the original "General" / unstyled appearance.
$dataRange.Style = "Normal"
